$d = $word.ActiveDocument

# Replace 1: "Mayornan/dunadónan di kuido por topa" -> "Mayornan/Edukadónan por topa"
$d.Content.Find.Execute("Mayornan/dunadónan di kuido por topa", $true, $true, $false, $false, $false, $true, 1, $false, "Mayornan/Edukadónan por topa", 2)

# Replace 2: "welanan/dunadónan di kuido," -> "welanan/edukadónan,"
$d.Content.Find.Execute("welanan/dunadónan di kuido,", $true, $true, $false, $false, $false, $true, 1, $false, "welanan/edukadónan,", 2)
